$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 107 / 109: Results column (E) changes from PASS to SKIP ---
$ws.Range("E107").Value = "SKIP"
$ws.Range("E109").Value = "SKIP"

# --- Add new test case rows 110-112 (new test case "Adding Test Case to search-OPQA-1242") ---
# Copy formatting (borders etc.) from an existing plain data row first, for all three rows.
$ws.Range("A105:E105").Copy($ws.Range("A110:E110"))
$ws.Range("A105:E105").Copy($ws.Range("A111:E111"))
$ws.Range("A105:E105").Copy($ws.Range("A112:E112"))

# Fill in cell values in the same order the strings were authored, so new shared-string
# entries come out in the expected sequence.
$ws.Range("B111").Value = "OPQA-581"
$ws.Range("C111").Value = "Verify that search results are sorted correctly by TIMES CITED field in SORT BY drop down in PATENTS search results page"
$ws.Range("B110").Value = "OPQA-592"
$ws.Range("C110").Value = "Verify that following filters are present in PATENTS search results page: a)Inventor b)IPC Codes c)Assignee"
$ws.Range("A110").Value = "TestCase_B109"
$ws.Range("A111").Value = "TestCase_B110"
$ws.Range("A112").Value = "TestCase_B111"
$ws.Range("B112").Value = "OPQA-1242"
$ws.Range("C112").Value = "Verify that more search results get displayed when user scrolls down in PEOPLE search results page."

$ws.Range("D110").Value = "Y"
$ws.Range("E110").Value = "SKIP"
$ws.Range("D111").Value = "Y"
$ws.Range("E111").Value = "SKIP"
$ws.Range("D112").Value = "Y"
$ws.Range("E112").Value = "PASS"

# --- Update sheet view: scroll position & active selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 98
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("D107").Select()
